$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ciArr = New-Object "object[,]" 24,7
$ciArr[0,0] = 0.1650118679366699
$ciArr[0,1] = 0.02393032774912029
$ciArr[0,2] = 0.2160984845063183
$ciArr[0,3] = 0.4256893306939133
$ciArr[0,4] = 0.2722042399704065
$ciArr[0,5] = 0.4471840647146337
$ciArr[0,6] = 0.3851305018112612
$ciArr[1,0] = 0.1679523250905532
$ciArr[1,1] = 0.02151435707205707
$ciArr[1,2] = 0.2007278981139464
$ciArr[1,3] = 0.4339240953785364
$ciArr[1,4] = 0.2810089727822245
$ciArr[1,5] = 0.4593102283487553
$ciArr[1,6] = 0.4036577243330548
$ciArr[2,0] = 0.1700221689005303
$ciArr[2,1] = 0.02002672195398247
$ciArr[2,2] = 0.1915133758166405
$ciArr[2,3] = 0.4397381122293709
$ciArr[2,4] = 0.2871200296797838
$ciArr[2,5] = 0.4672993803610055
$ciArr[2,6] = 0.4158589122269252
$ciArr[3,0] = 0.1709310360037222
$ciArr[3,1] = 0.01941949091180106
$ciArr[3,2] = 0.1878127011975437
$ciArr[3,3] = 0.4422949489621715
$ciArr[3,4] = 0.2897840617094971
$ciArr[3,5] = 0.4706902244931754
$ciArr[3,6] = 0.4210350861506331
$ciArr[4,0] = 0.1710858707823775
$ciArr[4,1] = 0.01931860147112019
$ciArr[4,2] = 0.1872014379605247
$ciArr[4,3] = 0.4427307518979475
$ciArr[4,4] = 0.290236812074447
$ciArr[4,5] = 0.4712613934504404
$ciArr[4,6] = 0.4219068116784612
$ciArr[5,0] = 0.1700341628590678
$ciArr[5,1] = 0.02001853662957842
$ciArr[5,2] = 0.1914632495558592
$ciArr[5,3] = 0.4397718390609597
$ciArr[5,4] = 0.2871552589781032
$ciArr[5,5] = 0.4673445650302597
$ciArr[5,6] = 0.4159278981145107
$ciArr[6,0] = 0.165970192328615
$ciArr[6,1] = 0.0230982063705838
$ciArr[6,2] = 0.2107512952723525
$ciArr[6,3] = 0.4283695870586683
$ciArr[6,4] = 0.2750916764502165
$ciArr[6,5] = 0.4512514300811006
$ciArr[6,6] = 0.3913452888872264
$ciArr[7,0] = 0.1601498267356476
$ciArr[7,1] = 0.02910217860512887
$ciArr[7,2] = 0.2504332322757534
$ciArr[7,3] = 0.41215756715156
$ciArr[7,4] = 0.2571873838706153
$ciArr[7,5] = 0.4240733686934135
$ciArr[7,6] = 0.3498433774851009
$ciArr[8,0] = 0.1572576118093849
$ciArr[8,1] = 0.0334902739569003
$ciArr[8,2] = 0.2808541466533967
$ciArr[8,3] = 0.4041811340156869
$ciArr[8,4] = 0.2477606118316018
$ciArr[8,5] = 0.4068679814776885
$ciArr[8,6] = 0.3236535072948996
$ciArr[9,0] = 0.1562593865574655
$ciArr[9,1] = 0.03548135900109628
$ciArr[9,2] = 0.2950003944079498
$ciArr[9,3] = 0.4014476501703896
$ciArr[9,4] = 0.244329945003102
$ciArr[9,5] = 0.3996606324285352
$ciArr[9,6] = 0.3127201361737271
$ciArr[10,0] = 0.1559285149356811
$ciArr[10,1] = 0.03623458449263239
$ciArr[10,2] = 0.3004042190729734
$ciArr[10,3] = 0.4005446782279023
$ciArr[10,4] = 0.2431582407832735
$ciArr[10,5] = 0.3970221729992289
$ciArr[10,6] = 0.3087249504288252
$ciArr[11,0] = 0.1559976533916227
$ciArr[11,1] = 0.03607239781581484
$ciArr[11,2] = 0.2992382757468022
$ciArr[11,3] = 0.4007332171730837
$ciArr[11,4] = 0.2434048554962231
$ciArr[11,5] = 0.3975863459572366
$ciArr[11,6] = 0.3095788690035093
$ciArr[12,0] = 0.1562312121681089
$ciArr[12,1] = 0.03554334256844527
$ciArr[12,2] = 0.2954440132812266
$ciArr[12,3] = 0.4013706939895769
$ciArr[12,4] = 0.2442309709113744
$ciArr[12,5] = 0.3994417348780956
$ciArr[12,6] = 0.3123885194803897
$ciArr[13,0] = 0.1563804571594574
$ciArr[13,1] = 0.03521918198563867
$ciArr[13,2] = 0.2931261153006091
$ciArr[13,3] = 0.4017784782482039
$ciArr[13,4] = 0.2447537063175815
$ciArr[13,5] = 0.4005900917009342
$ciArr[13,6] = 0.3141285220129717
$ciArr[14,0] = 0.157329361658654
$ciArr[14,1] = 0.03336004737616349
$ciArr[14,2] = 0.2799360742360903
$ciArr[14,3] = 0.404378076278384
$ciArr[14,4] = 0.2480024194066459
$ciArr[14,5] = 0.4073516111221664
$ciArr[14,6] = 0.3243881036530034
$ciArr[15,0] = 0.1579938281507083
$ciArr[15,1] = 0.03221820757835303
$ciArr[15,2] = 0.2719250773531314
$ciArr[15,3] = 0.4062045112240611
$ciArr[15,4] = 0.2502179433394147
$ciArr[15,5] = 0.4116594664898514
$ciArr[15,6] = 0.3309360164162811
$ciArr[16,0] = 0.1584057996221873
$ciArr[16,1] = 0.03156097429452132
$ciArr[16,2] = 0.2673461461739777
$ciArr[16,3] = 0.4073390646831498
$ciArr[16,4] = 0.2515727283353826
$ciArr[16,5] = 0.4141954333538393
$ciArr[16,6] = 0.3347942485634547
$ciArr[17,0] = 0.1585503580124907
$ciArr[17,1] = 0.03133836528142098
$ciArr[17,2] = 0.2658006622407072
$ciArr[17,3] = 0.4077375345990077
$ciArr[17,4] = 0.2520451333626568
$ciArr[17,5] = 0.4150640122341258
$ciArr[17,6] = 0.3361162666469681
$ciArr[18,0] = 0.1579200010779687
$ciArr[18,1] = 0.03233980796305502
$ciArr[18,2] = 0.2727748636712732
$ciArr[18,3] = 0.406001361966922
$ciArr[18,4] = 0.2499737388697625
$ciArr[18,5] = 0.4111948517701194
$ciArr[18,6] = 0.3302294244699731
$ciArr[19,0] = 0.1561613189118134
$ciArr[19,1] = 0.03569875957848012
$ciArr[19,2] = 0.2965571835421059
$ciArr[19,3] = 0.4011798377226015
$ciArr[19,4] = 0.2439848298076726
$ciArr[19,5] = 0.398894283782127
$ciArr[19,6] = 0.3115592879230853
$ciArr[20,0] = 0.1552873916136264
$ciArr[20,1] = 0.03788961992017903
$ciArr[20,2] = 0.312375331145148
$ciArr[20,3] = 0.3988005364314944
$ciArr[20,4] = 0.2408153566136946
$ciArr[20,5] = 0.3913852675825069
$ciArr[20,6] = 0.3002046405290066
$ciArr[21,0] = 0.1557281094453913
$ciArr[21,1] = 0.03672072635430368
$ciArr[21,2] = 0.3039067828815831
$ciArr[21,3] = 0.3999986406561078
$ciArr[21,4] = 0.2424374592740861
$ciArr[21,5] = 0.3953438869185817
$ciArr[21,6] = 0.3061859379690262
$ciArr[22,0] = 0.1579532850890786
$ciArr[22,1] = 0.03228483483901812
$ciArr[22,2] = 0.272390592322779
$ciArr[22,3] = 0.4060929427703073
$ciArr[22,4] = 0.2500838915642092
$ciArr[22,5] = 0.4114047194363266
$ciArr[22,6] = 0.3305485829843384
$ciArr[23,0] = 0.1614868861438623
$ciArr[23,1] = 0.02748190001256745
$ciArr[23,2] = 0.2394860096989646
$ciArr[23,3] = 0.4158660377126395
$ciArr[23,4] = 0.2613925287730012
$ciArr[23,5] = 0.430947359972663
$ciArr[23,6] = 0.3603310372728004
$ws.Range("C2:I25").Value = $ciArr

$lmArr = New-Object "object[,]" 24,2
$lmArr[0,0] = 0.4507150267031648
$lmArr[0,1] = 10.83166273276379
$lmArr[1,0] = 0.4056419359975507
$lmArr[1,1] = 9.493770058809446
$lmArr[2,0] = 0.3782603099191135
$lmArr[2,1] = 8.669003382941128
$lmArr[3,0] = 0.367172497386548
$lmArr[3,1] = 8.332052904883938
$lmArr[4,0] = 0.3653355346334308
$lmArr[4,1] = 8.276050537011486
$lmArr[5,0] = 0.3781104949182748
$lmArr[5,1] = 8.664462624953103
$lmArr[6,0] = 0.4351108030067792
$lmArr[6,1] = 10.3710257677198
$lmArr[7,0] = 0.5493872346675062
$lmArr[7,1] = 13.69286606534911
$lmArr[8,0] = 0.6351341822538359
$lmArr[8,1] = 16.12089203325957
$lmArr[9,0] = 0.6745949786964331
$lmArr[9,1] = 17.2234108905709
$lmArr[10,0] = 0.6896085990754557
$lmArr[10,1] = 17.64067521124707
$lmArr[11,0] = 0.6863719145287348
$lmArr[11,1] = 17.55081946101529
$lmArr[12,0] = 0.6758287097738958
$lmArr[12,1] = 17.25774388013184
$lmArr[13,0] = 0.6693800620402612
$lmArr[13,1] = 17.07819758574067
$lmArr[14,0] = 0.6325649149378592
$lmArr[14,1] = 16.04880388205987
$lmArr[15,0] = 0.6101000257180829
$lmArr[15,1] = 15.41682676331436
$lmArr[16,0] = 0.5972212117586366
$lmArr[16,1] = 15.05313829377025
$lmArr[17,0] = 0.5928677903539779
$lmArr[17,1] = 14.9299651590228
$lmArr[18,0] = 0.6124870259020838
$lmArr[18,1] = 15.4841213113101
$lmArr[19,0] = 0.6789235418462738
$lmArr[19,1] = 17.34383323297862
$lmArr[20,0] = 0.7227583431900939
$lmArr[20,1] = 18.55791729760517
$lmArr[21,0] = 0.6993230668615809
$lmArr[21,1] = 17.9100421117156
$lmArr[22,0] = 0.6114077496758057
$lmArr[22,1] = 15.45369856016708
$lmArr[23,0] = 0.5181787670616131
$lmArr[23,1] = 12.79662396037895
$ws.Range("L2:M25").Value = $lmArr

$oArr = New-Object "object[,]" 24,1
$oArr[0,0] = 1.360928806986848
$oArr[1,0] = 1.40452338796392
$oArr[2,0] = 1.433910271505596
$oArr[3,0] = 1.446533755224365
$oArr[4,0] = 1.448668708111555
$oArr[5,0] = 1.434077906127996
$oArr[6,0] = 1.375410099676998
$oArr[7,0] = 1.281626494062834
$oArr[8,0] = 1.226347551451823
$oArr[9,0] = 1.204302731451065
$oArr[10,0] = 1.196413179628678
$oArr[11,0] = 1.198091752636202
$oArr[12,0] = 1.203644396173587
$oArr[13,0] = 1.207105609243087
$oArr[14,0] = 1.227851696044141
$oArr[15,0] = 1.241381891427153
$oArr[16,0] = 1.249455314643328
$oArr[17,0] = 1.252238481732221
$oArr[18,0] = 1.239911351209145
$oArr[19,0] = 1.202000911739844
$oArr[20,0] = 1.179901690596296
$oArr[21,0] = 1.191447352502024
$oArr[22,0] = 1.240575265341008
$oArr[23,0] = 1.304653264975954
$ws.Range("O2:O25").Value = $oArr